$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Conc (ug/m3)"
$ws.Range("C1").Value = "Conc sci (ug/m3)"
